$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the A column labels for rows 3-7: the "[MG(FAn)-H2O+H]+" labels are replaced
# with "FAn_[MG-H2O+H]+" labels and re-ordered, while the "[M-(FAn)+H]+" rows shift up.
# New shared strings are appended in write order, so write FA3/FA2/FA1 in that order
# to match the target shared-strings table ordering.
$ws.Range("A3").Value = "FA2_[FA-H2O+H]+"
$ws.Range("A4").Value = "FA3_[FA-H2O+H]+"
$ws.Range("A7").Value = "FA3_[MG-H2O+H]+"
$ws.Range("A6").Value = "FA2_[MG-H2O+H]+"
$ws.Range("A5").Value = "FA1_[MG-H2O+H]+"

# Update the Group values (column C) for rows 5-7 from 1 to 3
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 3

# Update the active selection from C8 to C7
$ws.Range("C7").Select()
